$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column J values, keyed by row -> value
$values = @{
    4  = 2021
    5  = 1.5
    6  = 0.3
    7  = 0.8
    8  = 0.6
    9  = 1.8
    10 = 0.5
    11 = 0.8
    12 = 1.9
    13 = 4.4000000000000004
    14 = 0.4
}

foreach ($row in 4..14) {
    $srcCell = $ws.Range("I$row")
    $dstCell = $ws.Range("J$row")
    # Copy formatting/style from column I to the new column J cell
    $srcCell.Copy($dstCell)
    # Overwrite with the correct value for column J
    $dstCell.Value2 = $values[$row]
}

# Update the active selection to match the target state
$ws.Range("L10").Select()
